# Update "想去人数" (column F) values per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 453
$ws.Cells.Item(6, 6).Value = 597
$ws.Cells.Item(7, 6).Value = 62
$ws.Cells.Item(11, 6).Value = 364
$ws.Cells.Item(12, 6).Value = 1835
$ws.Cells.Item(13, 6).Value = 796
$ws.Cells.Item(15, 6).Value = 17
$ws.Cells.Item(16, 6).Value = 1558
$ws.Cells.Item(17, 6).Value = 1558
$ws.Cells.Item(18, 6).Value = 1296
$ws.Cells.Item(20, 6).Value = 1370
$ws.Cells.Item(21, 6).Value = 173
$ws.Cells.Item(22, 6).Value = 382
$ws.Cells.Item(25, 6).Value = 122
$ws.Cells.Item(26, 6).Value = 6806
$ws.Cells.Item(27, 6).Value = 7260
$ws.Cells.Item(28, 6).Value = 16
$ws.Cells.Item(32, 6).Value = 222
$ws.Cells.Item(37, 6).Value = 1336
$ws.Cells.Item(38, 6).Value = 203
$ws.Cells.Item(42, 6).Value = 1343
$ws.Cells.Item(43, 6).Value = 282
$ws.Cells.Item(47, 6).Value = 95

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(17, 6).Value = 2
$ws.Cells.Item(18, 6).Value = 258

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 229
$ws.Cells.Item(5, 6).Value = 94

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(4, 6).Value = 453
$ws.Cells.Item(6, 6).Value = 229
$ws.Cells.Item(7, 6).Value = 94
$ws.Cells.Item(9, 6).Value = 597
$ws.Cells.Item(10, 6).Value = 62
$ws.Cells.Item(15, 6).Value = 364
$ws.Cells.Item(16, 6).Value = 1835
$ws.Cells.Item(17, 6).Value = 796
$ws.Cells.Item(19, 6).Value = 17
$ws.Cells.Item(20, 6).Value = 1558
$ws.Cells.Item(21, 6).Value = 1558
$ws.Cells.Item(22, 6).Value = 1296
$ws.Cells.Item(24, 6).Value = 1370
$ws.Cells.Item(25, 6).Value = 173
$ws.Cells.Item(26, 6).Value = 382
$ws.Cells.Item(28, 6).Value = 122
$ws.Cells.Item(30, 6).Value = 6806
$ws.Cells.Item(31, 6).Value = 7260
$ws.Cells.Item(32, 6).Value = 222
$ws.Cells.Item(33, 6).Value = 1336
$ws.Cells.Item(34, 6).Value = 203
$ws.Cells.Item(43, 6).Value = 1343
$ws.Cells.Item(44, 6).Value = 282
$ws.Cells.Item(46, 6).Value = 95
$ws.Cells.Item(48, 6).Value = 2
$ws.Cells.Item(49, 6).Value = 258
